# Apply weekly cryptocurrency price/volume refresh (GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) and E (Volume(1h)) columns hold plain-text values (e.g. "23.532.42",
# "  +1.91%  "); force text format before writing so Excel does not silently
# reinterpret numeric-looking strings (like "131.30" or "0.2550") as numbers
# and drop significant trailing zeros / switch to scientific notation.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "23.532.42"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "1.639.64"
$ws.Range("E3").Value = "  +3.06%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "308.66"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "0.3774"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "52.89"
$ws.Range("E8").Value = "  +3.71%  "
$ws.Range("D9").Value = "0.3686"
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("D10").Value = "1.278"
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("D11").Value = "0.08216"
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("D12").Value = "0.9995"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "23.27"
$ws.Range("E13").Value = "  +4.09%  "
$ws.Range("D14").Value = "6.677"
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").Value = "0.00001284"
$ws.Range("E15").Value = "  +3.74%  "
$ws.Range("D16").Value = "7.483"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").Value = "1.641.60"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("D18").Value = "95.12"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("D19").Value = "0.06957"
$ws.Range("E19").Value = "  +2.84%  "
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").Value = "0.9982"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "23.527.94"
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("D24").Value = "12.96"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").Value = "3.096"
$ws.Range("E25").Value = "  +6.60%  "
$ws.Range("D26").Value = "2.419"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").Value = "21.42"
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("D28").Value = "151.35"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").Value = "5.323"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").Value = "136.07"
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").Value = "6.881"
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("D33").Value = "1.820.52"
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("D34").Value = "0.9852"
$ws.Range("E34").Value = "  +3.26%  "
$ws.Range("E35").Value = "  +5.27%  "
$ws.Range("D36").Value = "10.49"
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("D37").Value = "0.07491"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").Value = "6.236"
$ws.Range("E38").Value = "  +2.56%  "
$ws.Range("D39").Value = "0.2550"
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("D40").Value = "0.08883"
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("D41").Value = "1.399"
$ws.Range("E41").Value = "  +3.13%  "
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").Value = "12.66"
$ws.Range("E43").Value = "  +3.90%  "
$ws.Range("D44").Value = "16.21"
$ws.Range("E44").Value = "  +8.03%  "
$ws.Range("D45").Value = "0.6631"
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("D46").Value = "2.373"
$ws.Range("E46").Value = "  +4.11%  "
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "0.9987"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.08067"
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "131.30"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").Value = "1.226"
$ws.Range("E51").Value = "  +0.82%  "
